$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 loses its special ("Korean font") cell style -> back to default formatting.
# Do this before the copy/paste block below (interleaving a non-paste
# operation with an active clipboard/PasteSpecial sequence breaks the
# pending pastes).
$ws.Range("A2").ClearFormats()

# B1 / B2 ("BlackSmith_0" shop-id cells) now carry the same cell style as the
# rest of the row - copy that format from A1, which already uses it. Also
# stamp the same style onto the four new (still empty) NPC rows 8-11 in
# column A.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New column H for the existing 4 NPC dialogue rows - all 0.
$ws.Range("H1:H4").Value = 0

$null = $ws.Range("I13").Select()
